$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each year block (rows laid out as A/B/C/D quarters), the "B" quarter row
# and "C" quarter row need to swap places (all of columns A:E move together),
# while the "A" and "D" quarter rows stay where they are.
$pairs = @(
    @(3,4),
    @(7,8),
    @(11,12),
    @(15,16),
    @(19,20),
    @(23,24),
    @(27,28),
    @(31,32),
    @(35,36),
    @(39,40),
    @(43,44),
    @(47,48),
    @(51,52),
    @(55,56),
    @(59,60),
    @(63,64),
    @(67,68)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $a1 = $ws.Cells.Item($r1, 1).Value2
    $b1 = $ws.Cells.Item($r1, 2).Value2
    $c1 = $ws.Cells.Item($r1, 3).Value2
    $e1 = $ws.Cells.Item($r1, 5).Value2

    $a2 = $ws.Cells.Item($r2, 1).Value2
    $b2 = $ws.Cells.Item($r2, 2).Value2
    $c2 = $ws.Cells.Item($r2, 3).Value2
    $e2 = $ws.Cells.Item($r2, 5).Value2

    $ws.Cells.Item($r1, 1).Value2 = $a2
    $ws.Cells.Item($r1, 2).Value2 = $b2
    $ws.Cells.Item($r1, 3).Value2 = $c2
    $ws.Cells.Item($r1, 5).Value2 = $e2

    $ws.Cells.Item($r2, 1).Value2 = $a1
    $ws.Cells.Item($r2, 2).Value2 = $b1
    $ws.Cells.Item($r2, 3).Value2 = $c1
    $ws.Cells.Item($r2, 5).Value2 = $e1

    # column D only carries real data for the 2016+ blocks; swap it too when present
    $d1 = $ws.Cells.Item($r1, 4).Value2
    $d2 = $ws.Cells.Item($r2, 4).Value2
    if (($d1 -ne "") -or ($d2 -ne "")) {
        $ws.Cells.Item($r1, 4).Value2 = $d2
        $ws.Cells.Item($r2, 4).Value2 = $d1
    }
}

# Drop the "粗钢产销率" (F) and "粗钢销售量" (G) columns entirely.
$ws.Range("F1:G1").EntireColumn.Delete()
